# Added second reference diet and ran results.
# Update the "type" column (F) for the "Spices" food group rows
# (rows 1731-1947) from "B" to "A" to reflect the second reference diet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1731:F1947").Value = "A"
